$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "IsRanged" header in H6, matching the header style used by E6:G6 (s="3")
$ws.Range("H6").Value = "IsRanged"
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New boolean column (IsRanged) for the data rows
$ws.Range("H7").Value = $false
$ws.Range("H8").Value = $false
$ws.Range("H9").Value = $false
$ws.Range("H10").Value = $false
$ws.Range("H11").Value = $true

# Update C9 (Bob Jones Health) value from 100 to 50
$ws.Range("C9").Value = 50

# Update the current selection to E6:E11 (active cell E6)
$ws.Range("E6:E11").Select()
